$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: replace values with 2-decimal "custom accuracy" rounded values
$row5Values = @{
    "B5"  = 6.3
    "C5"  = 4.43
    "D5"  = 0.07000000000000001
    "E5"  = 11.88
    "F5"  = 10.19
    "G5"  = 4.82
    "H5"  = 19.2
    "I5"  = 6.6
    "J5"  = 3.13
    "K5"  = 5.08
    "L5"  = 4.89
    "M5"  = 4.94
    "N5"  = 1.45
    "O5"  = 4.27
    "P5"  = 6.73
    "Q5"  = 3.59
    "R5"  = 0.28
    "S5"  = 0.09
    "T5"  = 61.16
    "U5"  = 12.67
    "V5"  = 4.33
    "W5"  = 8.65
    "X5"  = 4.32
    "Y5"  = 0.59
    "Z5"  = 9.23
    "AA5" = 3.66
    "AB5" = 3.06
    "AC5" = 3.65
    "AD5" = 5.61
    "AE5" = 0.52
    "AF5" = 17.25
    "AG5" = 2.22
    "AH5" = 5.04
}

foreach ($addr in $row5Values.Keys) {
    $ws.Range($addr).Value = $row5Values[$addr]
}

# Row 6 is removed entirely from the used range
$ws.Range("A6:AH6").EntireRow.Delete()
